$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Change text_field to email_field for user email field*") {
        $p.Range.Delete()
        break
    }
}
